# Apply the "commit" edits:
# 1. Rename Sheet3 -> Log, make it the active tab with tabSelected, move topLeft/selection
# 2. Clear H3 / H10 "Update LoadContents function" text on DDAS sheet
# 3. Set H13 on DDAS sheet to "Need to implement search function"
# 4. Update DDAS sheet view: remove tabSelected, move topLeftCell to B2, selection to B13
# 5. Populate Log (Sheet3) sheet with header rows ("Working copy path:" etc and
#    "Files Updated"/"Files Created"/"Context updated"/"Updated On"/"Created On")
# 6. Widen column B on Log sheet to fit the new text

$wb = $excel.ActiveWorkbook

# --- Sheet: DDAS ---
$ddas = $wb.Worksheets.Item("DDAS")

# Clear the two "Update LoadContents function" notes; they become blank and
# pick up the same (default) formatting as the other blank H-column cells.
$ddas.Range("H3").Value = $null
$ddas.Range("H10").Value = $null
$ddas.Range("H2").Copy()
$ddas.Range("H3").PasteSpecial(-4122)
$ddas.Range("H10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rename Sheet3 -> Log ---
$log = $wb.Worksheets.Item("Sheet3")
$log.Name = "Log"

# Populate the Log sheet headers (order matters for shared-string allocation).
$log.Range("A1").Value = "Working copy path:"
$log.Range("B1").Value = "C:\Development\p926-ddas"

$log.Range("A3").Value = "Files Updated"
$log.Range("D3").Value = "Files Created"
$log.Range("B3").Value = "Context updated"
$log.Range("C3").Value = "Updated On"
$log.Range("E3").Value = "Created On"

# New note on row 13 of DDAS (added last so it gets the final shared-string slot).
$ddas.Range("H13").Value = "Need to implement search function"

$log.Columns.Item(2).ColumnWidth = 26.140625

# Update the selection / view state for the DDAS sheet.
$ddas.Activate()
$ddas.Range("B13").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2

# Make the Log sheet the active/selected tab.
$log.Activate()
$log.Range("A1").Select()

$wb.Save()
